# Updated cryptos list on Wed Oct 18 15:29:45 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    # Force the written value to be stored as text (matches the source
    # data's inline-string cells), even when it looks numeric (e.g.
    # "211.83"), without leaving a residual NumberFormat on the cell.
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

# Row 2 - Bitcoin
Set-TextValue "D2" "28.325.16"
$ws.Range("E2").Value = "  -0.82%  "

# Row 3 - Ethereum
Set-TextValue "D3" "1.575.29"
$ws.Range("E3").Value = "  -0.08%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.23%  "

# Row 5 - BNB
Set-TextValue "D5" "211.83"
$ws.Range("E5").Value = "  -0.51%  "

# Row 6 - XRP
$ws.Range("E6").Value = "  -0.62%  "

# Row 7 - USDC
$ws.Range("E7").Value = "  +0.25%  "

# Row 8 - OKB
Set-TextValue "D8" "44.48"
$ws.Range("E8").Value = "  -3.74%  "

# Row 9 - Solana
Set-TextValue "D9" "23.83"
$ws.Range("E9").Value = "  -1.30%  "

# Row 10 - Cardano
$ws.Range("E10").Value = "  -0.92%  "

# Row 11 - Dogecoin
$ws.Range("E11").Value = "  -0.90%  "

# Row 12 - TRON
$ws.Range("E12").Value = "  +1.61%  "

# Row 13 - WrappedliquidstakedEther2.0
Set-TextValue "D13" "1.800.95"
$ws.Range("E13").Value = "  +0.00%  "

# Row 14 - WrappedEther
Set-TextValue "D14" "1.575.64"
$ws.Range("E14").Value = "  -0.03%  "

# Row 15 - Polkadot
Set-TextValue "D15" "3.68"
$ws.Range("E15").Value = "  -0.81%  "

# Row 16 - Polygon
$ws.Range("E16").Value = "  -1.23%  "

# Row 17 - WrappedBTC
Set-TextValue "D17" "28.351.15"
$ws.Range("E17").Value = "  -0.65%  "

# Row 18 - Litecoin
Set-TextValue "D18" "61.57"
$ws.Range("E18").Value = "  -1.40%  "

# Row 19 - BitcoinCash
Set-TextValue "D19" "229.95"
$ws.Range("E19").Value = "  -0.06%  "

# Row 20 - Chainlink
$ws.Range("E20").Value = "  +0.14%  "

# Row 21 - ShibaInu
Set-TextValue "D21" "0.0₃0685"
$ws.Range("E21").Value = "  -1.25%  "

# Row 22 - Dai
$ws.Range("E22").Value = "  +0.24%  "

# Row 23 - Uniswap
$ws.Range("E23").Value = "  +0.11%  "

# Row 24 - Avalanche
Set-TextValue "D24" "9.03"
$ws.Range("E24").Value = "  -1.53%  "

# Row 25 - Toncoin
$ws.Range("E25").Value = "  +2.18%  "

# Row 26 - Monero
Set-TextValue "D26" "151.91"
$ws.Range("E26").Value = "  +0.40%  "

# Row 27 - EthereumClassic
Set-TextValue "D27" "14.95"
$ws.Range("E27").Value = "  -0.65%  "

# Row 28 - Cosmos
Set-TextValue "D28" "6.36"
$ws.Range("E28").Value = "  -1.63%  "

# Row 29 - Stellar
$ws.Range("E29").Value = "  -1.56%  "

# Row 30 - BinanceUSD
$ws.Range("E30").Value = "  +0.28%  "

# Row 31 - Hedera
Set-TextValue "D31" "0.0480"
$ws.Range("E31").Value = "  +3.35%  "

# Row 32 - PancakeSwap
$ws.Range("E32").Value = "  -3.99%  "

# Row 33 - Filecoin
Set-TextValue "D33" "3.19"
$ws.Range("E33").Value = "  -0.54%  "

# Row 34 - InternetComputer(DFINITY)
Set-TextValue "D34" "3.06"
$ws.Range("E34").Value = "  -2.09%  "

# Row 35 - Maker
Set-TextValue "D35" "1.386.73"
$ws.Range("E35").Value = "  -0.52%  "

# Row 36 - TrustWalletToken
$ws.Range("E36").Value = "  +5.80%  "

# Row 37 - LidoDAOToken
$ws.Range("E37").Value = "  -3.15%  "

# Row 38 - HuobiToken: unchanged

# Row 39 - MXToken
$ws.Range("E39").Value = "  +2.93%  "

# Row 40 - VeChain
$ws.Range("E40").Value = "  -1.88%  "

# Row 41 - ImmutableX
$ws.Range("E41").Value = "  -2.52%  "

# Row 42 - PaxDollar
$ws.Range("E42").Value = "  +0.25%  "

# Row 43 - RenderToken
$ws.Range("E43").Value = "  +1.58%  "

# Row 44 - ARBITRUM
$ws.Range("E44").Value = "  -1.21%  "

# Row 45 - Kaspa
Set-TextValue "D45" "0.0461"
$ws.Range("E45").Value = "  -0.25%  "

# Row 46 - FraxShare
$ws.Range("E46").Value = "  -4.31%  "

# Row 47 & 48 - Aave and WEMIXToken swapped places (with updated prices)
$ws.Range("B47").Value = "WEMIXToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
Set-TextValue "D47" "0.924"
$ws.Range("E47").Value = "  -5.52%  "

$ws.Range("B48").Value = "Aave"
$ws.Range("C48").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
Set-TextValue "D48" "62.26"
$ws.Range("E48").Value = "  -0.58%  "

# Row 49 - RocketPoolETH
Set-TextValue "D49" "1.712.47"

# Row 50 - mCoin
$ws.Range("E50").Value = "  +0.73%  "

# Row 51 - Quant
Set-TextValue "D51" "85.44"
$ws.Range("E51").Value = "  -0.22%  "
